$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Osm"
$ws.Range("C2").Value = "Osmr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.025991
$ws.Range("H2").Value = 0.051982
$ws.Range("I2").Value = 0.0002490198180052299
$ws.Range("J2").Value = 0.0001660269933407085
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 62.245413
$ws.Range("N2").Value = 124.490826
$ws.Range("O2").Value = 0.3724206900939814
$ws.Range("P2").Value = 0.3142670423823846
$ws.Range("Q2").Value = 1.617820529283
$ws.Range("R2").Value = 6.471282117132
$ws.Range("S2").Value = 0.00009274013246858537
$ws.Range("T2").Value = 0.0000521768121528243

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Osm"
$ws.Range("C3").Value = "Osmr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.025991
$ws.Range("H3").Value = 0.051982
$ws.Range("I3").Value = 0.0002490198180052299
$ws.Range("J3").Value = 0.0001660269933407085
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 61.55916833333333
$ws.Range("N3").Value = 184.677505
$ws.Range("O3").Value = 0.3683148178695765
$ws.Range("P3").Value = 0.4662034557543063
$ws.Range("Q3").Value = 1.599984344151667
$ws.Range("R3").Value = 9.59990606491
$ws.Range("S3").Value = 0.00009171768891451132
$ws.Range("T3").Value = 0.00007740235804393549

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Osm"
$ws.Range("C4").Value = "Osmr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.025991
$ws.Range("H4").Value = 0.051982
$ws.Range("I4").Value = 0.0002490198180052299
$ws.Range("J4").Value = 0.0001660269933407085
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 43.0359755
$ws.Range("N4").Value = 86.071951
$ws.Range("O4").Value = 0.2574886553420037
$ws.Range("P4").Value = 0.2172816932940226
$ws.Range("Q4").Value = 1.1185480392205
$ws.Range("R4").Value = 4.474192156882
$ws.Range("S4").Value = 0.00006411977809167712
$ws.Range("T4").Value = 0.00003607462624558454

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Osm"
$ws.Range("C5").Value = "Osmr"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.025991
$ws.Range("H5").Value = 0.051982
$ws.Range("I5").Value = 0.0002490198180052299
$ws.Range("J5").Value = 0.0001660269933407085
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.159175
$ws.Range("N5").Value = 0.477525
$ws.Range("O5").Value = 0.0009523603505644583
$ws.Range("P5").Value = 0.001205473320689897
$ws.Range("Q5").Value = 0.004137117425
$ws.Range("R5").Value = 0.02482270455
$ws.Range("S5").Value = 0.0000002371566011729583
$ws.Range("T5").Value = 0.0000002001411109865833

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Osm"
$ws.Range("C6").Value = "Osmr"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.025991
$ws.Range("H6").Value = 0.051982
$ws.Range("I6").Value = 0.0002490198180052299
$ws.Range("J6").Value = 0.0001660269933407085
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1376336666666667
$ws.Range("N6").Value = 0.412901
$ws.Range("O6").Value = 0.0008234763438739657
$ws.Range("P6").Value = 0.001042335248596784
$ws.Range("Q6").Value = 0.003577236630333334
$ws.Range("R6").Value = 0.021463419782
$ws.Range("S6").Value = 0.000000205061929283107
$ws.Range("T6").Value = 0.0000001730557873775639

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Osm"
$ws.Range("C7").Value = "Osmr"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 45.64166666666667
$ws.Range("H7").Value = 136.925
$ws.Range("I7").Value = 0.4372928908771768
$ws.Range("J7").Value = 0.4373291920891176
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 62.245413
$ws.Range("N7").Value = 124.490826
$ws.Range("O7").Value = 0.3724206900939814
$ws.Range("P7").Value = 0.3142670423823846
$ws.Range("Q7").Value = 2840.984391675001
$ws.Range("R7").Value = 17045.90635005
$ws.Range("S7").Value = 0.1628569201936703
$ws.Range("T7").Value = 0.1374381517453247

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Osm"
$ws.Range("C8").Value = "Osmr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 45.64166666666667
$ws.Range("H8").Value = 136.925
$ws.Range("I8").Value = 0.4372928908771768
$ws.Range("J8").Value = 0.4373291920891176
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 61.55916833333333
$ws.Range("N8").Value = 184.677505
$ws.Range("O8").Value = 0.3683148178695765
$ws.Range("P8").Value = 0.4662034557543063
$ws.Range("Q8").Value = 2809.663041347223
$ws.Range("R8").Value = 25286.967372125
$ws.Range("S8").Value = 0.161061451459088
$ws.Range("T8").Value = 0.2038843806541854

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Osm"
$ws.Range("C9").Value = "Osmr"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 45.64166666666667
$ws.Range("H9").Value = 136.925
$ws.Range("I9").Value = 0.4372928908771768
$ws.Range("J9").Value = 0.4373291920891176
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 43.0359755
$ws.Range("N9").Value = 86.071951
$ws.Range("O9").Value = 0.2574886553420037
$ws.Range("P9").Value = 0.2172816932940226
$ws.Range("Q9").Value = 1964.233648445834
$ws.Range("R9").Value = 11785.401890675
$ws.Range("S9").Value = 0.1125979584625818
$ws.Range("T9").Value = 0.09502362738403032

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Osm"
$ws.Range("C10").Value = "Osmr"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 45.64166666666667
$ws.Range("H10").Value = 136.925
$ws.Range("I10").Value = 0.4372928908771768
$ws.Range("J10").Value = 0.4373291920891176
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.159175
$ws.Range("N10").Value = 0.477525
$ws.Range("O10").Value = 0.0009523603505644583
$ws.Range("P10").Value = 0.001205473320689897
$ws.Range("Q10").Value = 7.265012291666667
$ws.Range("R10").Value = 65.385110625
$ws.Range("S10").Value = 0.0004164604108551335
$ws.Range("T10").Value = 0.0005271886734222986

# Row 11
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Osm"
$ws.Range("C11").Value = "Osmr"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 45.64166666666667
$ws.Range("H11").Value = 136.925
$ws.Range("I11").Value = 0.4372928908771768
$ws.Range("J11").Value = 0.4373291920891176
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.1376336666666667
$ws.Range("N11").Value = 0.412901
$ws.Range("O11").Value = 0.0008234763438739657
$ws.Range("P11").Value = 0.001042335248596784
$ws.Range("Q11").Value = 6.281829936111112
$ws.Range("R11").Value = 56.53646942500001
$ws.Range("S11").Value = 0.0003601003509816146
$ws.Range("T11").Value = 0.0004558436321548411

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Osm"
$ws.Range("C12").Value = "Osmr"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 38.34823466666666
$ws.Range("H12").Value = 115.044704
$ws.Range("I12").Value = 0.3674145056948628
$ws.Range("J12").Value = 0.367445006057708
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 62.245413
$ws.Range("N12").Value = 124.490826
$ws.Range("O12").Value = 0.3724206900939814
$ws.Range("P12").Value = 0.3142670423823846
$ws.Range("Q12").Value = 2387.001704647584
$ws.Range("R12").Value = 14322.0102278855
$ws.Range("S12").Value = 0.1368327637614199
$ws.Range("T12").Value = 0.1154758552919333

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Osm"
$ws.Range("C13").Value = "Osmr"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 38.34823466666666
$ws.Range("H13").Value = 115.044704
$ws.Range("I13").Value = 0.3674145056948628
$ws.Range("J13").Value = 0.367445006057708
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 61.55916833333333
$ws.Range("N13").Value = 184.677505
$ws.Range("O13").Value = 0.3683148178695765
$ws.Range("P13").Value = 0.4662034557543063
$ws.Range("Q13").Value = 2360.685433131502
$ws.Range("R13").Value = 21246.16889818352
$ws.Range("S13").Value = 0.1353242067476439
$ws.Range("T13").Value = 0.1713041316237655

# Row 14
$ws.Range("A14").Value = "Neutrophils"
$ws.Range("B14").Value = "Osm"
$ws.Range("C14").Value = "Osmr"
$ws.Range("D14").Value = "MuSCs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 38.34823466666666
$ws.Range("H14").Value = 115.044704
$ws.Range("I14").Value = 0.3674145056948628
$ws.Range("J14").Value = 0.367445006057708
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 43.0359755
$ws.Range("N14").Value = 86.071951
$ws.Range("O14").Value = 0.2574886553420037
$ws.Range("P14").Value = 0.2172816932940226
$ws.Range("Q14").Value = 1650.353687582917
$ws.Range("R14").Value = 9902.122125497504
$ws.Range("S14").Value = 0.09460506702451718
$ws.Range("T14").Value = 0.07983907310865117

# Row 15
$ws.Range("A15").Value = "Neutrophils"
$ws.Range("B15").Value = "Osm"
$ws.Range("C15").Value = "Osmr"
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 38.34823466666666
$ws.Range("H15").Value = 115.044704
$ws.Range("I15").Value = 0.3674145056948628
$ws.Range("J15").Value = 0.367445006057708
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.159175
$ws.Range("N15").Value = 0.477525
$ws.Range("O15").Value = 0.0009523603505644583
$ws.Range("P15").Value = 0.001205473320689897
$ws.Range("Q15").Value = 6.104080253066665
$ws.Range("R15").Value = 54.9367222776
$ws.Range("S15").Value = 0.0003499110074460267
$ws.Range("T15").Value = 0.0004429451516233047

# Row 16
$ws.Range("A16").Value = "Neutrophils"
$ws.Range("B16").Value = "Osm"
$ws.Range("C16").Value = "Osmr"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 38.34823466666666
$ws.Range("H16").Value = 115.044704
$ws.Range("I16").Value = 0.3674145056948628
$ws.Range("J16").Value = 0.367445006057708
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.1376336666666667
$ws.Range("N16").Value = 0.412901
$ws.Range("O16").Value = 0.0008234763438739657
$ws.Range("P16").Value = 0.001042335248596784
$ws.Range("Q16").Value = 5.278008147367111
$ws.Range("R16").Value = 47.502073326304
$ws.Range("S16").Value = 0.000302557153835866
$ws.Range("T16").Value = 0.0003830008817348078

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Osm"
$ws.Range("C17").Value = "Osmr"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 20.35732666666667
$ws.Range("H17").Value = 61.07198
$ws.Range("I17").Value = 0.1950435836099552
$ws.Range("J17").Value = 0.1950597748598338
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 62.245413
$ws.Range("N17").Value = 124.490826
$ws.Range("O17").Value = 0.3724206900939814
$ws.Range("P17").Value = 0.3142670423823846
$ws.Range("Q17").Value = 1267.15020594258
$ws.Range("R17").Value = 7602.901235655479
$ws.Range("S17").Value = 0.07263826600642267
$ws.Range("T17").Value = 0.06130085853297378

# Row 18
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Osm"
$ws.Range("C18").Value = "Osmr"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 20.35732666666667
$ws.Range("H18").Value = 61.07198
$ws.Range("I18").Value = 0.1950435836099552
$ws.Range("J18").Value = 0.1950597748598338
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 61.55916833333333
$ws.Range("N18").Value = 184.677505
$ws.Range("O18").Value = 0.3683148178695765
$ws.Range("P18").Value = 0.4662034557543063
$ws.Range("Q18").Value = 1253.180099089989
$ws.Range("R18").Value = 11278.6208918099
$ws.Range("S18").Value = 0.07183744197393016
$ws.Range("T18").Value = 0.09093754111831147

# Row 19
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Osm"
$ws.Range("C19").Value = "Osmr"
$ws.Range("D19").Value = "MuSCs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 20.35732666666667
$ws.Range("H19").Value = 61.07198
$ws.Range("I19").Value = 0.1950435836099552
$ws.Range("J19").Value = 0.1950597748598338
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 43.0359755
$ws.Range("N19").Value = 86.071951
$ws.Range("O19").Value = 0.2574886553420037
$ws.Range("P19").Value = 0.2172816932940226
$ws.Range("Q19").Value = 876.0974116721633
$ws.Range("R19").Value = 5256.58447003298
$ws.Range("S19").Value = 0.05022151007681304
$ws.Range("T19").Value = 0.0423829181750955

# Row 20
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Osm"
$ws.Range("C20").Value = "Osmr"
$ws.Range("D20").Value = "Neutrophils"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 20.35732666666667
$ws.Range("H20").Value = 61.07198
$ws.Range("I20").Value = 0.1950435836099552
$ws.Range("J20").Value = 0.1950597748598338
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.159175
$ws.Range("N20").Value = 0.477525
$ws.Range("O20").Value = 0.0009523603505644583
$ws.Range("P20").Value = 0.001205473320689897
$ws.Range("Q20").Value = 3.240377472166666
$ws.Range("R20").Value = 29.1633972495
$ws.Range("S20").Value = 0.0001857517756621252
$ws.Range("T20").Value = 0.0002351393545333076

# Row 21
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Osm"
$ws.Range("C21").Value = "Osmr"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 20.35732666666667
$ws.Range("H21").Value = 61.07198
$ws.Range("I21").Value = 0.1950435836099552
$ws.Range("J21").Value = 0.1950597748598338
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.1376336666666667
$ws.Range("N21").Value = 0.412901
$ws.Range("O21").Value = 0.0008234763438739657
$ws.Range("P21").Value = 0.001042335248596784
$ws.Range("Q21").Value = 2.801853512664445
$ws.Range("R21").Value = 25.21668161398
$ws.Range("S21").Value = 0.000160613777127202
$ws.Range("T21").Value = 0.0002033176789197576
